$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated expectation/variance values for NVDA and AMAT (row 2 & 3)
$ws.Range("B2").Value = 0.004791757665818062
$ws.Range("C2").Value = 0.0007813032928621058

$ws.Range("B3").Value = 0.002288414088021248
$ws.Range("C3").Value = 0.0003758880376423061

# Row 6 previously held TXN; replace it with the CTAS figures (previously row 7)
$ws.Range("A6").Value = "CTAS"
$ws.Range("B6").Value = 0.00128081856973723
$ws.Range("C6").Value = 0.0001526322847375969

# Remove the now-duplicate CTAS row (old row 7) entirely, shrinking the used range
$ws.Range("A7:C7").Delete()
